$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy the existing header style (bold, bordered, centered) from H1
# to the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF)
$dataI = @(7, 8, 7, 9, 8, 9, 6, 7, 9, 6, 8, 7, 10, 7, 1, 7, 6, 7, 6, 6, 9, 5)
$dataJ = @(7, 8, 7, 9, 9, 9, 6, 8, 9, 7, 8, 7, 10, 7, 2, 7, 6, 7, 6, 6, 9, 5)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
